$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MZ_Reference")

$data = New-Object 'object[,]' 16,9
$data[0,0] = "cas"
$data[0,1] = "name"
$data[0,2] = "short_display_name"
$data[0,3] = "id"
$data[0,4] = "trt"
$data[0,5] = "mz0"
$data[0,6] = "mz1"
$data[0,7] = "mz2"
$data[0,8] = "mz3"
$data[1,0] = "2093-28-9"
$data[1,1] = "2,6-Dichlorophenyl-4'-nitrophenyl ether"
$data[1,2] = "2,6-DCP-4’-NPE"
$data[1,3] = "CP2458"
$data[1,4] = 10.4388
$data[1,5] = 282.9797
$data[1,6] = 139.0545
$data[1,7] = 204.015
$data[1,8] = 254.9796
$data[2,0] = "91-59-8"
$data[2,1] = "2-Naphthylamine"
$data[2,2] = "2-Naphthylamine"
$data[2,3] = "CP2535"
$data[2,4] = 7.364
$data[2,5] = 142.0731
$data[2,6] = 117.0574
$data[2,7] = 116.0495
$data[2,8] = 107.0732
$data[3,0] = "91-59-8"
$data[3,1] = "2-Naphthylamine"
$data[3,2] = "2-Naphthylamine"
$data[3,3] = "CP3014"
$data[3,4] = 7.3405
$data[3,5] = 143.073
$data[3,6] = 117.0574
$data[3,7] = 118.0652
$data[3,8] = 119.073
$data[4,0] = "92-67-1"
$data[4,1] = "4-Aminobiphenyl"
$data[4,2] = "4-ABP"
$data[4,3] = "CP2518"
$data[4,4] = 6.6385
$data[4,5] = 169.0884
$data[4,6] = 168.0808
$data[4,7] = 154.0653
$data[4,8] = 141.0699
$data[5,0] = "92-67-1"
$data[5,1] = "4-aminobiphenyl"
$data[5,2] = "4-ABP"
$data[5,3] = "CP3002"
$data[5,4] = 5.366
$data[5,5] = 169.0888
$data[5,6] = 168.081
$data[5,7] = 167.0731
$data[5,8] = 170.0924
$data[6,0] = "92-87-5"
$data[6,1] = "Benzidine"
$data[6,2] = "Benzidine"
$data[6,3] = "CP2215"
$data[6,4] = 9.6825
$data[6,5] = 185.0805
$data[6,6] = 183.1334
$data[6,7] = 91.0542
$data[6,8] = 91.0542
$data[7,0] = "92-87-5"
$data[7,1] = "Benzidine"
$data[7,2] = "Benzidine"
$data[7,3] = "CP3094"
$data[7,4] = 5.6984
$data[7,5] = 184.0996
$data[7,6] = 182.0838
$data[7,7] = 183.0872
$data[7,8] = 181.0761
$data[8,0] = "50-32-8"
$data[8,1] = "Benz[a]pyrene"
$data[8,2] = "Benzo[a]pyrene"
$data[8,3] = "CP2221"
$data[8,4] = 16.9112
$data[8,5] = 252.0942
$data[8,6] = 250.0786
$data[8,7] = 253.0974
$data[8,8] = 248.0629
$data[9,0] = "50-32-8"
$data[9,1] = "Benz[a]pyrene"
$data[9,2] = "Benzo[a]pyrene"
$data[9,3] = "CP3028"
$data[9,4] = 16.6163
$data[9,5] = 252.0937
$data[9,6] = 250.0781
$data[9,7] = 253.0969
$data[9,8] = 126.0463
$data[10,0] = "68359-37-5"
$data[10,1] = "Cyfluthrin"
$data[10,2] = "Cyfluthrin"
$data[10,3] = "CP3153"
$data[10,4] = 8.3914
$data[10,5] = 433.0677
$data[10,6] = 163.0542
$data[10,7] = 165.0699
$data[10,8] = 167.0855
$data[11,0] = "117-84-0"
$data[11,1] = "Di-n-octyl phthalate"
$data[11,2] = "DNOP"
$data[11,3] = "CP2187"
$data[11,4] = 17.847
$data[11,5] = 149.1238
$data[11,6] = 150.1359
$data[11,7] = 275.1057
$data[11,8] = 151.1392
$data[12,0] = "117-84-0"
$data[12,1] = "Di-n-octyl phthalate"
$data[12,2] = "DNOP"
$data[12,3] = "CP3056"
$data[12,4] = 14.1913
$data[12,5] = 149.0238
$data[12,6] = 279.159
$data[12,7] = 150.0302
$data[12,8] = 167.0339
$data[13,0] = "2642-71-9"
$data[13,1] = "Azinphos ethyl"
$data[13,2] = "Guthion"
$data[13,3] = "CP2002"
$data[13,4] = 15.1056
$data[13,5] = 132.0445
$data[13,6] = 104.0495
$data[13,7] = 137.0056
$data[13,8] = 133.0477
$data[14,0] = "95-53-4"
$data[14,1] = "o-Toluidine"
$data[14,2] = "o-Toluidine"
$data[14,3] = "CP2551"
$data[14,4] = 5.797
$data[14,5] = 106.0653
$data[14,6] = 107.073
$data[14,7] = 108.0683
$data[14,8] = 105.0573
$data[15,0] = "95-53-4"
$data[15,1] = "o-Toluidine"
$data[15,2] = "o-Toluidine"
$data[15,3] = "CP3017"
$data[15,4] = 7.8706
$data[15,5] = 107.073
$data[15,6] = 106.0651
$data[15,7] = 108.0764
$data[15,8] = 105.0336

$ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(16,9)).Value = $data

# Clear old rows 17-18 that are no longer part of the table
$ws.Range("A17:I18").ClearContents()
